$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# Cells that carry the build timestamp: A2 & A6 on "About", S2:S17 on
# "Boundaries and methane sources". Update each in place by replacing the
# old timestamp with the new one so the rest of each string is preserved
# exactly.

$targets = @(
    @{ Sheet = "About"; Cell = "A2" },
    @{ Sheet = "About"; Cell = "A6" }
)

for ($row = 2; $row -le 17; $row++) {
    $targets += @{ Sheet = "Boundaries and methane sources"; Cell = "S" + $row }
}

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $cell = $ws.Range($t.Cell)
    $current = $cell.Value()
    $updated = $current.Replace($oldStamp, $newStamp)
    $cell.Value = $updated
}
